$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 / Row 23 swap: Chainlink <-> Toncoin ---
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"

# --- Price (column D) updates ---
# Cells whose new text could be misread as a number need the column
# format forced to Text ("@") first so Excel keeps them as strings,
# matching the workbook's original inline-string price formatting.
$ws.Range("D2").Value = "26.740.59"
$ws.Range("D3").Value = "1.639.00"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D12").Value = "1.867.65"
$ws.Range("D13").Value = "1.634.18"
$ws.Range("D17").Value = "26.730.89"
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "211.41"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.34"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.19"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D34").Value = "1.275.60"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.811"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.803"
$ws.Range("D43").Value = "1.777.42"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.80"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.15"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.57"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  -0.07%  "
